$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.88789165019989
$ws.Range("B1").Value = 2.820541620254517
$ws.Range("C1").Value = 3.116150617599487
$ws.Range("D1").Value = 2.670375823974609
$ws.Range("E1").Value = 1.048498272895813
